$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update B52 text: append new sentence about storage full behavior
$ws.Range("B52").Value = "加入根据飞行速度调节翅膀频率。加入落地`n物品消失action。仓库满时，改为把第一个挤掉"

# Update C52 text: expand bug note, enable wrap text, and bump row height
$ws.Range("C52").Value = "水晶球削球好像有bug。未找出。看代码没看出来`n仓库削球效果没有层次感"
$ws.Range("C52").WrapText = $true

# Update D52 value from 2 to 4
$ws.Range("D52").Value = 4

# Row height for row 52 grows because of the extra wrapped line in C52
$ws.Rows.Item(52).RowHeight = 40.5

# Update the active selection on the sheet view
$ws.Activate()
$ws.Range("C55").Select()
